# Updated symbol list on Wed Dec 28 04:49:04 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) quotes for the crypto-price snapshot sheet,
# plus a handful of coin identity / label corrections (rows 18, 41, 43) that
# came through in this run's pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a textual (non-numeric) value, same as the
    # inlineStr cells already used throughout this sheet, then restore the
    # cell's original (default) style so nothing else about it changes.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Column D price refresh -------------------------------------------------
Set-TextValue $ws.Range("D2")  "245.28"
Set-TextValue $ws.Range("D3")  "23.77"
Set-TextValue $ws.Range("D4")  "5.325"
Set-TextValue $ws.Range("D5")  "0.05786"
Set-TextValue $ws.Range("D6")  "6.469"
Set-TextValue $ws.Range("D7")  "3.335"
Set-TextValue $ws.Range("D8")  "0.8124"
Set-TextValue $ws.Range("D9")  "0.8908"
Set-TextValue $ws.Range("D10") "0.1391"
Set-TextValue $ws.Range("D11") "0.07351"
Set-TextValue $ws.Range("D12") "0.03090"
Set-TextValue $ws.Range("D13") "0.03062"
Set-TextValue $ws.Range("D14") "0.09360"
Set-TextValue $ws.Range("D15") "3.849"
Set-TextValue $ws.Range("D16") "0.001544"
Set-TextValue $ws.Range("D17") "0.04715"

# Row 18 (One / ONE): price refresh + "Worstin24h" tag appended to label
Set-TextValue $ws.Range("D18") "0.0006004"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue $ws.Range("D19") "0.005911"
Set-TextValue $ws.Range("D20") "0.001294"
Set-TextValue $ws.Range("D22") "0.00008801"
Set-TextValue $ws.Range("D25") "0.3178"
Set-TextValue $ws.Range("D40") "0.03809"

# Row 41: BKEXToken -> KickToken (coin identity swap with row 43)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006341"
$ws.Range("E41").Value = "40KickTokenKICK"

Set-TextValue $ws.Range("D42") "0.004101"

# Row 43: KickToken -> BKEXToken (coin identity swap with row 41)
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D43") "0.1055"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue $ws.Range("D44") "0.007839"
Set-TextValue $ws.Range("D45") "0.00005473"
Set-TextValue $ws.Range("D48") "0.001844"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("D50") "0.0002000"
